$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TeamStats")

# Add the missing "average velocity" figure to the overall summary row,
# matching the other AVERAGE() formulas already present in row 18.
$ws.Range("B18").Formula = "=AVERAGE(B2:B16)"

# Reflect the cell the author was last working in when the results were
# committed.
$ws.Range("B18").Select()
